# edit.ps1 - apply the "step1 a new database doesn't work" change set
# Strategy: locate target text with Find.Execute to obtain a Range, then
# use Range.InsertXML with a minimal flat-OPC wrapper to inject exact
# run/paragraph level OOXML (runs with w:lang, w:proofErr marks,
# bookmarks, lastRenderedPageBreak, and brand-new paragraphs) without
# the engine silently re-merging same-formatted runs.

$d = $word.ActiveDocument

function New-OpenXmlPackage([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphRuns($range, [string]$runsXml) {
    # Replace the run content of a single paragraph (everything except
    # its trailing paragraph mark), keeping the paragraph's own pPr.
    $xml = New-OpenXmlPackage("<w:p>$runsXml</w:p>")
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) "Как менять первое меню" paragraph -> split into several runs and
#    change the wording to "Как менять Название первого меню"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Как менять первое меню")
$target = $d.Range($r.Start, $r.End)
$runs =
  '<w:r><w:t>Как менять</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>Название</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> перво</w:t></w:r>' +
  '<w:r><w:t>го</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> меню</w:t></w:r>'
Set-ParagraphRuns $target $runs

# ---------------------------------------------------------------------
# 2) + 3) Insert two brand-new paragraphs right after it:
#    "Где хранится данные БД при работе из debug_handlers.py?"
#    "Как можно обратиться к данной БД из Редактора?"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Как менять Название первого меню")
$insertPoint = $d.Range($r.End, $r.End)
$newParas =
  '<w:p><w:r><w:t xml:space="preserve">Где хранится данные БД при работе из </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>debug</w:t></w:r>' +
  '<w:r><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>handlers</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>py</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>?</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>Как можно обратиться к данной БД из Редактора?</w:t></w:r></w:p>'
$insertPoint.InsertXML((New-OpenXmlPackage $newParas))

# ---------------------------------------------------------------------
# 4) Duplicate the blank paragraph right before "Главное меню приложения"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Главное меню приложения")
$para = $d.Paragraphs.Item($d.Range(0, $r.Start).Paragraphs.Count)
$blankBefore = $para.Previous()
$insertPoint = $d.Range($blankBefore.Range.End, $blankBefore.Range.End)
$insertPoint.InsertXML((New-OpenXmlPackage '<w:p/>'))

# ---------------------------------------------------------------------
# 5) Add lastRenderedPageBreak before "Главное меню приложения"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Главное меню приложения")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full '<w:r><w:lastRenderedPageBreak/><w:t>Главное меню приложения</w:t></w:r>'

# ---------------------------------------------------------------------
# 6) Remove lastRenderedPageBreak before "Отображается: " (the first one,
#    which precedes the "Кнопка «Закупки»" list)
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Отображается: ")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full '<w:r><w:t xml:space="preserve">Отображается: </w:t></w:r>'

# ---------------------------------------------------------------------
# 7) Split "Кнопка «Закупки»" with a bookmarkStart between « and Закупки»
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Кнопка «Закупки»")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full (
  '<w:r><w:t>Кнопка «</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_Hlk115779802"/>' +
  '<w:r><w:t>Закупки»</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 8) Split "Кнопка «Товары»" with a bookmarkEnd between Товары and »
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Кнопка «Товары»")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full (
  '<w:r><w:t>Кнопка «Товары</w:t></w:r>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t>»</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 9) Add lastRenderedPageBreak before "1.1 Список документов закупки"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("1.1 Список документов закупки")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full '<w:r><w:lastRenderedPageBreak/><w:t>1.1 Список документов закупки</w:t></w:r>'

# ---------------------------------------------------------------------
# 10) Insert the "Обработчики ..." paragraph plus a blank paragraph right
#     before "Доступные действия:" (the one following the "Закупка" table)
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Таблица со списком документов «Закупка»")
$insertPoint = $d.Range($r.End, $r.End)
$handlersPara =
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:t>Обработчики</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>ПриЗапуске</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> / </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>Ввода</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>listbuy</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>on</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>start</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> / </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>listbuy</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>on</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>input</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$blankPara = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$insertPoint.InsertXML((New-OpenXmlPackage ($handlersPara + $blankPara)))

# ---------------------------------------------------------------------
# 11) Remove lastRenderedPageBreak before "Если при открытии экрана
#     список документов пуст ..."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Если при открытии экрана список документов пуст – автоматически создается новый документ.")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full '<w:r><w:t>Если при открытии экрана список документов пуст – автоматически создается новый документ.</w:t></w:r>'

# ---------------------------------------------------------------------
# 12) Add lastRenderedPageBreak before "1.1.1 Документ закупки"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("1.1.1 Документ закупки")
$full = $d.Range($r.Start, $r.End)
Set-ParagraphRuns $full '<w:r><w:lastRenderedPageBreak/><w:t>1.1.1 Документ закупки</w:t></w:r>'

Write-Output "done"
